$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "YBuqd553"
$ws.Range("B2").Value = 23072720
$ws.Range("C2").Value = "peglstf18"
$ws.Range("D2").Value = "Bv%&J93u"
$ws.Range("F2").Value = "JoyaOucZ"
$ws.Range("G2").Value = "rCEY"
